$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 19 reflects more work done and a later pending-date:
# - ActlHours (P19) increased from 9 to 12
# - Weekday of actual delivery (R19) moved from Friday to Monday
# - ActDeliveryDate comment (Q19) date pushed from 2020-02-14 to 2020-02-17
$ws.Range("P19").Value = 12
$ws.Range("R19").Value = "Monday"
$ws.Range("Q19").Value = "2020-02-17 Pending"

$wb.Save()
